# Trade #56 closed at 2026-02-17 13:29:47 - unknown UNKNOWN +0.000%
#
# Updates the roll-up metrics on "Summary" and "Strategy Status" to account
# for the newly-closed MarketMaking trade, and appends the trade's row to
# both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# --- Summary sheet: refresh aggregate stats -------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.54   # Current Capital
$summary.Range("B4").Value = -2.45     # Total P&L $
$summary.Range("B5").Value = -0.88     # Total P&L %
$summary.Range("B6").Value = 56        # Total Trades
$summary.Range("B7").Value = 23        # Winning Trades
$summary.Range("B9").Value = 41.07     # Win Rate %

# --- Strategy Status sheet: refresh the MarketMaking row ------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.54000000000001  # Capital
$status.Range("D4").Value = 56                 # Trades
$status.Range("E4").Value = -2.45              # P&L $
$status.Range("F4").Value = -2.46              # P&L %
$status.Range("G4").Value = 41.07              # Win Rate %

# --- Append the new trade row to a trade log sheet -------------------------
function Add-TradeRow {
    param($ws, $r)

    $ws.Cells.Item($r, 1).Value = 56

    # Date/Time columns look like dates to Excel's auto-detection, so force
    # them to text (matching every other row in these logs) and then drop
    # back to the workbook's default "Normal" style so no stray per-cell
    # number format sticks around once the literal text is locked in.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = "2026-02-17"
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = "13:29:41"
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).Value = "MarketMaking"
    $ws.Cells.Item($r, 5).Value = "UP"
    $ws.Cells.Item($r, 6).Value = 0.9399999999999999
    $ws.Cells.Item($r, 7).Value = 0.98
    $ws.Cells.Item($r, 8).Value = "CLOSED"
    $ws.Cells.Item($r, 9).Value = 4.2553
    $ws.Cells.Item($r, 10).Value = 0.04
    $ws.Cells.Item($r, 11).Value = 97.54000000000001
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0.6
    $ws.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($r, 16).Value = "early_exit"
    $ws.Cells.Item($r, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 57

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 57
